$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value = -7.702
$ws.Range("D4").Value = -7.879
$ws.Range("C7").Value = -13.497
$ws.Range("B8").Value = 6.962000000000001
$ws.Range("B10").Value = 6.017
$ws.Range("D11").Value = -7.225
$ws.Range("B12").Value = 5.624000000000001
$ws.Range("C14").Value = -13.094
$ws.Range("D14").Value = -7.806
$ws.Range("C15").Value = -13.745
$ws.Range("B18").Value = 5.544
$ws.Range("C18").Value = -13.415
$ws.Range("D18").Value = -8.606999999999999
$ws.Range("D19").Value = -8.255000000000001
$ws.Range("C20").Value = -12.684
$ws.Range("D21").Value = -8.219999999999999
$ws.Range("B25").Value = 6.751
$ws.Range("D27").Value = -8.568999999999999
$ws.Range("C29").Value = -12.491
$ws.Range("C30").Value = -12.338
$ws.Range("C31").Value = -13.397
$ws.Range("D31").Value = -8.389999999999999
$ws.Range("C35").Value = -12.518
$ws.Range("B37").Value = 8.416
$ws.Range("D38").Value = -7.869
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.238
$ws.Range("C44").Value = -12.395
$ws.Range("D44").Value = -7.672999999999999
$ws.Range("D47").Value = -7.486
$ws.Range("C50").Value = -13.542
$ws.Range("C54").Value = -12.72
$ws.Range("B55").Value = 5.151
$ws.Range("D56").Value = -8.440999999999999
$ws.Range("D58").Value = -8.134
$ws.Range("D65").Value = -7.633999999999999
$ws.Range("B68").Value = 5.403
$ws.Range("C68").Value = -11.083
$ws.Range("D73").Value = -8.318000000000001
$ws.Range("C76").Value = -13.46
$ws.Range("B77").Value = 5.476999999999999
$ws.Range("B78").Value = 7.114
$ws.Range("B79").Value = 5.17
$ws.Range("B80").Value = 8.273999999999999
$ws.Range("B81").Value = 5.575
$ws.Range("B82").Value = 5.867
$ws.Range("B84").Value = 6.145
$ws.Range("C87").Value = -13.215
$ws.Range("C88").Value = -12.768
$ws.Range("D90").Value = -7.556
$ws.Range("C92").Value = -11.372
$ws.Range("D92").Value = -6.488
$ws.Range("D94").Value = -6.798999999999999
$ws.Range("D95").Value = -7.833000000000001
$ws.Range("C96").Value = -12.705
$ws.Range("C98").Value = -13.45
$ws.Range("B101").Value = 8.975
$ws.Range("C101").Value = -13.022
$ws.Range("D101").Value = -7.783999999999999
$ws.Range("B102").Value = 7.468000000000001
$ws.Range("C102").Value = -12.906
